$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("Data")
$wsEntropy = $wb.Worksheets.Item("Entropy")

# ---------------------------------------------------------------------------
# "Data" sheet: finish entering the last three experiment subjects (rows
# 9-11) for the columns that were still blank (D, I, N, S, X, AA, AB, AC).
# Row 12 holds shared AVERAGE-style formulas (=SUM(..)/10) that will
# recalculate automatically once these inputs are present.
# ---------------------------------------------------------------------------

# Row 9 (subject 8)
$wsData.Range("D9").Value = 93
$wsData.Range("I9").Value = 0.57889999999999997
$wsData.Range("N9").Value = 5449.2502999999997
$wsData.Range("S9").Value = 0.20899999999999999
$wsData.Range("X9").Value = 107.1152
$wsData.Range("AA9").Value = 0.090899999999999995
$wsData.Range("AB9").Value = 0
$wsData.Range("AC9").Value = 1

# Row 10 (subject 9)
$wsData.Range("D10").Value = 68
$wsData.Range("I10").Value = 0.90349999999999997
$wsData.Range("N10").Value = 3799.5549999999998
$wsData.Range("S10").Value = 0.16400000000000001
$wsData.Range("X10").Value = 140.5213
$wsData.Range("AA10").Value = 0
$wsData.Range("AB10").Value = 0
$wsData.Range("AC10").Value = 1

# Row 11 (subject 10)
$wsData.Range("D11").Value = 135
$wsData.Range("I11").Value = 0.86160000000000003
$wsData.Range("N11").Value = 3186.0774999999999
$wsData.Range("S11").Value = 0.26440000000000002
$wsData.Range("X11").Value = 175.23570000000001
$wsData.Range("AA11").Value = 0
$wsData.Range("AB11").Value = 0
$wsData.Range("AC11").Value = 1

# ---------------------------------------------------------------------------
# "Entropy" sheet: the summary table (rows 3-5, one row per algorithm) now
# picks up the refreshed "Hybrid" row (row 5) from the completed Data-sheet
# averages, plus the two QM values (column H) that feed rows 3 & 4.
# ---------------------------------------------------------------------------
$wsEntropy.Range("H3").Value = 0.01363
$wsEntropy.Range("H4").Value = 0.0001

$wsEntropy.Range("C5").Value = 96.9
$wsEntropy.Range("D5").Value = 0.74052000000000007
$wsEntropy.Range("E5").Value = 4606.7493699999995
$wsEntropy.Range("F5").Value = 0.20927000000000001
$wsEntropy.Range("G5").Value = 133.47863999999998
$wsEntropy.Range("H5").Value = 1

$wsEntropy.Range("AH8").Value = 1.0912750733522381

# ---------------------------------------------------------------------------
# View state: the saved file now has "Entropy" as the active tab, with a
# fresh selection on each sheet.
# ---------------------------------------------------------------------------
$wsData.Range("B12:AC12").Select()
$wsEntropy.Activate()
$wsEntropy.Range("H3").Select()
